$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RestrictJudgments")

# --- Row 3 (DESCRIPTION): update EN/KO text, drop the other language columns
#     (the old per-language "Kills/Mata/Tue/Lam ban/..." strings are removed).
$ws.Range("D3:H3").ClearContents()
$ws.Range("B3").Value = "Punishes player upon specified judgments."
$ws.Range("C3").Value = "특정 판정에서 플레이어를 제한합니다."

# --- Insert a new row 4 for the RESTRICT_HEADER key (EN/KO only)
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "RESTRICT_HEADER"
$ws.Range("B4").Value = "Restrict following judgments:"
$ws.Range("C4").Value = "제한할 판정 목록:"

# (old row 4 "RESTRICT" is now row 5, shifted down automatically, content unchanged)

# --- Insert two new rows before the CUSTOM_DEATH row (currently row 6)
#     for CUSTOM_HEADER and RESTRICT_ACTION (EN/KO only)
$ws.Range("A6:A7").EntireRow.Insert()

$ws.Range("A6").Value = "CUSTOM_HEADER"
$ws.Range("B6").Value = "Restrict Action:"
$ws.Range("C6").Value = "플레이어 제한 행동:"

$ws.Range("A7").Value = "RESTRICT_ACTION"
$ws.Range("B7").Value = "Restriction Method:"
$ws.Range("C7").Value = "제한 방식:"

# (old row 5 "CUSTOM_DEATH" is now row 8, shifted down automatically, content unchanged)
